# Oracle SQL cheatsheet - Table of contents update
# 0.4 Constraints still WIP. Next step: CHECK Constraint

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above row 7 for the new "VERIFY..." Functions entry ---
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "VERIFY, LENGTH, COUNTW, FIND, REPLACE, etc."
$ws.Range("B7").Value = "Functions"

# --- 2. Unhide every data row (no rows remain hidden in the final layout) ---
$ws.Rows.Item(2).Hidden = $false
$ws.Rows.Item(3).Hidden = $false
$ws.Rows.Item(13).Hidden = $false
$ws.Rows.Item(14).Hidden = $false
$ws.Rows.Item(15).Hidden = $false
$ws.Rows.Item(16).Hidden = $false
$ws.Rows.Item(17).Hidden = $false

# --- 3. ALTER TABLE (now row 12) is marked as Done ---
$ws.Range("C12").Value = 1

# --- 4. Append the new "ADD / DROP CONSTRAINT" row (row 18), still WIP ---
$ws.Range("A18").Value = "ADD / DROP CONSTRAINT"
$ws.Range("B18").Value = "DB management utilities"
$ws.Range("C18").Value = "WIP"

# --- 5. Refresh the AutoFilter to cover the grown range and clear any filter criteria ---
$ws.AutoFilterMode = $false
$ws.Range("A1:C18").AutoFilter()

# --- 6. Restore the on-screen view: scrolled to row 7, C7 selected ---
$ws.Range("C7").Select()
$excel.ActiveWindow.ScrollRow = 7

Write-Host "edit applied"
